$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old row 17 (English "Syllabus:" row with the long English program text),
# which shifts subsequent rows up by one and yields the new dimension A1:C24.
$ws.Rows.Item(17).Delete()

# Row height fix-ups for the rows whose content/size changed as a result of the
# upstream data regeneration (rows 13 and 15 gain the heights that used to belong
# to the rows that shifted into their place).
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120

# Cell content fix-ups to match the final published sheet exactly.
$ws.Range('B10').Value = '230696 - Carlos José Todero Peixoto'
$ws.Range('C10').Value = '230696 - Carlos José Todero Peixoto'
$ws.Range('A13').Value = 'Programa resumido:'
$ws.Range('B13').Value = 'Semestral'
$ws.Range('C13').Value = 'Semestral'
$ws.Range('A14').Value = 'Short syllabus:'
$ws.Range('B14').Value = 'Geometrical Optics. Introduction to Modern Physics: wave nature of matter, relativity and introduction to quantum mechanics.'
$ws.Range('C14').Value = 'Geometrical Optics. Introduction to Modern Physics: wave nature of matter, relativity and introduction to quantum mechanics.'
$ws.Range('A15').Value = 'Programa:'
$ws.Range('B15').Value = '01/01/2018'
$ws.Range('C15').Value = '01/01/2018'
$ws.Range('A16').Value = 'Syllabus:'
$ws.Range('B16').Value = '1) Geometrical Optics: basic concepts.2) Interference: Young''s experience; coherence; interference figures; the Michelson interferometer.3) Diffraction.4) Polarization.5) Relativity: the postulates of relativity, Lorentz transformations, simultaneity, time and length; linear momentum, work and energy;6) Early days of quantum theory: the hypothesis of Planck; the photoelectric effect, quantization of the photon; De Broglie waves, the Compton effect, the electron diffraction, interference;7) Basic principles of quantum mechanics: the uncertainty principle; the Schrödinger equation.'
$ws.Range('C16').Value = '1) Geometrical Optics: basic concepts.2) Interference: Young''s experience; coherence; interference figures; the Michelson interferometer.3) Diffraction.4) Polarization.5) Relativity: the postulates of relativity, Lorentz transformations, simultaneity, time and length; linear momentum, work and energy;6) Early days of quantum theory: the hypothesis of Planck; the photoelectric effect, quantization of the photon; De Broglie waves, the Compton effect, the electron diffraction, interference;7) Basic principles of quantum mechanics: the uncertainty principle; the Schrödinger equation.'
$ws.Range('B18').Value = '230696 - Carlos José Todero Peixoto'
$ws.Range('C18').Value = '230696 - Carlos José Todero Peixoto'
$ws.Range('B19').Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range('C19').Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range('B20').Value = 'NF≥ 5,0.'
$ws.Range('C20').Value = 'NF≥ 5,0.'
$ws.Range('B21').Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range('C21').Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
"Edit applied"
